$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Find-ParaByText($txt) {
    $result = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -eq $txt) {
            $result = $p
        }
    }
    return $result
}

# Insert raw OOXML (one or more <w:p>...</w:p> fragments) immediately AFTER
# the given paragraph, without disturbing the paragraph itself or its
# neighbours. We anchor on a point strictly inside the paragraph's text
# (Start+1) because anchoring exactly at a paragraph boundary causes the
# host to splice into / devour the adjacent paragraph.
function Insert-XmlAfterParagraph($para, [string]$xmlFragment) {
    $anchor = $para.Range.Start + 1
    $r = $d.Range($anchor, $anchor)
    $r.InsertXML($xmlFragment)
}

# Replace an existing paragraph (identified by its exact, current text
# including trailing CR) with a brand new paragraph described by raw OOXML.
# We do this by inserting the replacement right after the original and then
# deleting the whole original paragraph (this also removes any bookmarks
# that lived inside it).
function Replace-Paragraph([string]$oldText, [string]$xmlFragment) {
    $target = Find-ParaByText $oldText
    Insert-XmlAfterParagraph $target $xmlFragment
    $old = Find-ParaByText $oldText
    $old.Range.Delete()
}

# ---------------------------------------------------------------------
# 1. Turn the "Power consumption" paragraph into an empty paragraph, and
#    push its text down into a brand-new paragraph that follows 4 more
#    brand-new blank paragraphs (5 blank paragraphs total end up between
#    the old content and the relocated "Power consumption" text).
# ---------------------------------------------------------------------
$emptyP = "<w:p $wns><w:pPr><w:pStyle w:val=`"NoSpacing`"/></w:pPr></w:p>"
$powerParaXml = "<w:p $wns><w:pPr><w:pStyle w:val=`"NoSpacing`"/></w:pPr><w:r><w:t>Power consumption</w:t></w:r></w:p>"
$newContent = $emptyP + $emptyP + $emptyP + $emptyP + $powerParaXml

$powerPara = Find-ParaByText "Power consumption`r"
# Remember the original paragraph's own (pre-insert) span: since we only
# ever insert content strictly *after* it, this span stays valid.
$origStart = $powerPara.Range.Start
$origEnd = $powerPara.Range.End
Insert-XmlAfterParagraph $powerPara $newContent

# Clear the original paragraph's own text (but keep the paragraph / its
# original rsid attributes), leaving it as the first of the 5 blanks.
$textOnly = $d.Range($origStart, $origEnd - 1)
$textOnly.Text = ""

# ---------------------------------------------------------------------
# 2. Replace "-sonar : 15mA x 4 x5V" paragraph
# ---------------------------------------------------------------------
$sonarXml = "<w:p $wns><w:pPr><w:pStyle w:val=`"NoSpacing`"/></w:pPr>" +
  "<w:r><w:t xml:space=`"preserve`">-sonar : 15mA x </w:t></w:r>" +
  "<w:r><w:t>////////////</w:t></w:r>" +
  "<w:r><w:t>4 x5V</w:t></w:r>" +
  "<w:r><w:t xml:space=`"preserve`"> = 0.3A</w:t></w:r>" +
  "</w:p>"
Replace-Paragraph "-sonar : 15mA x 4 x5V`r" $sonarXml

# ---------------------------------------------------------------------
# 3. Replace "-infrared: 0.33mA x 5V" paragraph
# ---------------------------------------------------------------------
$infraredXml = "<w:p $wns><w:pPr><w:pStyle w:val=`"NoSpacing`"/></w:pPr>" +
  "<w:r><w:t xml:space=`"preserve`">-infrared: 0.33mA x </w:t></w:r>" +
  "<w:r><w:t>/////////////</w:t></w:r>" +
  "<w:r><w:t>5V</w:t></w:r>" +
  "<w:r><w:t xml:space=`"preserve`"> = </w:t></w:r>" +
  "<w:r><w:t>0.00165A</w:t></w:r>" +
  "</w:p>"
Replace-Paragraph "-infrared: 0.33mA x 5V`r" $infraredXml

# ---------------------------------------------------------------------
# 4. Replace "-altimu-10 v4: 6mA x 3.3V " paragraph
# ---------------------------------------------------------------------
$altimuXml = "<w:p $wns><w:pPr><w:pStyle w:val=`"NoSpacing`"/></w:pPr>" +
  "<w:r><w:t>-altimu-10 v4: 6mA x</w:t></w:r>" +
  "<w:r><w:t>//////////////</w:t></w:r>" +
  "<w:r><w:t xml:space=`"preserve`"> 3.3V </w:t></w:r>" +
  "<w:r><w:t xml:space=`"preserve`">= </w:t></w:r>" +
  "<w:r><w:t>0.0198A</w:t></w:r>" +
  "</w:p>"
Replace-Paragraph "-altimu-10 v4: 6mA x 3.3V `r" $altimuXml

# ---------------------------------------------------------------------
# 5. Replace "-motor: 80mA x2 x3.3V" paragraph (this also removes the
#    _GoBack bookmark that used to live at the end of this paragraph)
# ---------------------------------------------------------------------
$motorXml = "<w:p $wns><w:pPr><w:pStyle w:val=`"NoSpacing`"/></w:pPr>" +
  "<w:r><w:t>-motor: 80mA x2 x</w:t></w:r>" +
  "<w:r><w:t>/////////////////</w:t></w:r>" +
  "<w:r><w:t>3.3V</w:t></w:r>" +
  "<w:r><w:t>=</w:t></w:r>" +
  "</w:p>"
Replace-Paragraph "-motor: 80mA x2 x3.3V`r" $motorXml

# ---------------------------------------------------------------------
# 6. Append new paragraphs after the (now bookmark-less) motor paragraph:
#    Rpi bareboard, Arduino, blank, Total, blank, AA eneloop (with bookmark)
# ---------------------------------------------------------------------
$rpiRpr = "<w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/>" +
  "<w:color w:val=`"222222`"/><w:shd w:val=`"clear`" w:color=`"auto`" w:fill=`"F2F1F0`"/></w:rPr>"
$rpiXml = "<w:p $wns><w:pPr><w:pStyle w:val=`"NoSpacing`"/>$rpiRpr</w:pPr>" +
  "<w:r><w:t>-Rpi</w:t></w:r>" +
  "<w:r><w:t xml:space=`"preserve`"> bareboard </w:t></w:r>" +
  "<w:r><w:t xml:space=`"preserve`">: </w:t></w:r>" +
  "<w:r>$rpiRpr<w:t>200mA</w:t></w:r>" +
  "<w:r>$rpiRpr<w:t xml:space=`"preserve`"> x </w:t></w:r>" +
  "<w:r>$rpiRpr<w:t>////////////////</w:t></w:r>" +
  "<w:r>$rpiRpr<w:t>3</w:t></w:r>" +
  "<w:r>$rpiRpr<w:t>.3</w:t></w:r>" +
  "<w:r>$rpiRpr<w:t>V</w:t></w:r>" +
  "</w:p>"

$arduinoRpr = "<w:rPr><w:rFonts w:ascii=`"TyponineSans Regular 18`" w:hAnsi=`"TyponineSans Regular 18`"/>" +
  "<w:color w:val=`"222222`"/><w:shd w:val=`"clear`" w:color=`"auto`" w:fill=`"F5F5F5`"/></w:rPr>"
$arduinoXml = "<w:p $wns><w:pPr><w:pStyle w:val=`"NoSpacing`"/>$arduinoRpr</w:pPr>" +
  "<w:r>$arduinoRpr<w:t xml:space=`"preserve`">-Arduino - </w:t></w:r>" +
  "<w:r>$arduinoRpr<w:t>200.0 mA</w:t></w:r>" +
  "<w:r>$arduinoRpr<w:t xml:space=`"preserve`"> x </w:t></w:r>" +
  "<w:r>$arduinoRpr<w:t>///////////////////</w:t></w:r>" +
  "<w:r>$arduinoRpr<w:t>5V</w:t></w:r>" +
  "</w:p>"

$blankXml = $emptyP

$totalXml = "<w:p $wns><w:pPr><w:pStyle w:val=`"NoSpacing`"/></w:pPr>" +
  "<w:r><w:t xml:space=`"preserve`">Total: </w:t></w:r>" +
  "<w:r><w:t xml:space=`"preserve`">15m + 0.33m + 6m + 80m + 200m </w:t></w:r>" +
  "<w:r><w:t xml:space=`"preserve`">+ </w:t></w:r>" +
  "<w:r><w:t xml:space=`"preserve`">200m = </w:t></w:r>" +
  "<w:r><w:t>501.33mA</w:t></w:r>" +
  "</w:p>"

$eneloopXml = "<w:p $wns><w:pPr><w:pStyle w:val=`"NoSpacing`"/></w:pPr>" +
  "<w:r><w:t>AA eneloop</w:t></w:r>" +
  "<w:r><w:t xml:space=`"preserve`"> recharagable: </w:t></w:r>" +
  "<w:r><w:t xml:space=`"preserve`"> 200mAh</w:t></w:r>" +
  "<w:r><w:t xml:space=`"preserve`"> X 6  </w:t></w:r>" +
  "<w:r><w:sym w:font=`"Wingdings`" w:char=`"F0E8`"/></w:r>" +
  "<w:r><w:t xml:space=`"preserve`"> last for 2.4hours</w:t></w:r>" +
  "<w:bookmarkStart w:id=`"4`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"4`"/>" +
  "</w:p>"

$allNewTail = $rpiXml + $arduinoXml + $blankXml + $totalXml + $blankXml + $eneloopXml

$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)
Insert-XmlAfterParagraph $lastPara $allNewTail

Write-Host "Edit complete."
